# This script reorders the observation records on the active sheet so that
# rows 2,4,5 and rows 6,7,9,10,11 each take on the record that used to live
# in a different row of the same group (the data was "re-sorted" upstream).
# Row positions stay the same; only the cell contents belonging to each
# record move between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: read a set of single-letter/double-letter column values from a row
# into a hashtable, using Value() (parentheses needed to force evaluation).
function Get-RowVals($ws, $row, $cols) {
    $d = @{}
    foreach ($c in $cols) {
        $d[$c] = $ws.Range($c + $row).Value()
    }
    return $d
}

# Helper: write a hashtable of column -> value back onto a row.
function Set-RowVals($ws, $row, $d) {
    foreach ($c in $d.Keys) {
        $ws.Range($c + $row).Value = $d[$c]
    }
}

# Columns that carry the actual record data (species/observation fields).
$recordCols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# --- Group 1: rows 2, 4, 5 rotate among themselves -------------------------
# new row2 = old row4, new row4 = old row5, new row5 = old row2
$orig2 = Get-RowVals $ws 2 $recordCols
$orig4 = Get-RowVals $ws 4 $recordCols
$orig5 = Get-RowVals $ws 5 $recordCols

Set-RowVals $ws 2 $orig4
Set-RowVals $ws 4 $orig5
Set-RowVals $ws 5 $orig2

# --- Group 2: rows 6, 7, 9, 10, 11 rotate among themselves ------------------
# new row6 = old row10, new row7 = old row6, new row9 = old row11,
# new row10 = old row9, new row11 = old row7
$orig6  = Get-RowVals $ws 6  $recordCols
$orig7  = Get-RowVals $ws 7  $recordCols
$orig9  = Get-RowVals $ws 9  $recordCols
$orig10 = Get-RowVals $ws 10 $recordCols
$orig11 = Get-RowVals $ws 11 $recordCols

Set-RowVals $ws 6  $orig10
Set-RowVals $ws 7  $orig6
Set-RowVals $ws 9  $orig11
Set-RowVals $ws 10 $orig9
Set-RowVals $ws 11 $orig7

# Startdatum/Slutdatum (Y, AA) only actually change value between rows 9 and
# 10 (the other rows in the group share the same date). These look like
# dates, so a plain .Value assignment would get auto-parsed into a date
# serial number instead of staying plain text; use copy/paste of the values
# instead, through a scratch cell, to keep them as literal text.
$ws.Range("Y9").Copy()
$ws.Range("BZ1").PasteSpecial(-4163)
$ws.Range("AA9").Copy()
$ws.Range("BZ2").PasteSpecial(-4163)

$ws.Range("Y10").Copy()
$ws.Range("Y9").PasteSpecial(-4163)
$ws.Range("AA10").Copy()
$ws.Range("AA9").PasteSpecial(-4163)

$ws.Range("BZ1").Copy()
$ws.Range("Y10").PasteSpecial(-4163)
$ws.Range("BZ2").Copy()
$ws.Range("AA10").PasteSpecial(-4163)

$ws.Range("BZ1:BZ2").ClearContents()

# The (empty) "Bestämningsmetod" placeholder cell moves from row 6 to row 7,
# and from row 9 to row 10, as part of the same record rotation. Row 6 and
# row 9 lose it (clear to blank); rows 7/10 end up blank as well, which is
# already their rotated-in state.
$ws.Range("AF6").ClearContents()
$ws.Range("AF9").ClearContents()
